$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.005990743637085
$ws.Range("B1").Value = 3.462850570678711
$ws.Range("C1").Value = 3.076259613037109
$ws.Range("D1").Value = 3.298145771026611
$ws.Range("E1").Value = 1.008774995803833
